$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menu")
$ws.Activate()
